# Jan 21 - Cart Lang input
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for rows 2-24 (columns A-E): Environment, Brand, Campaign, Categories, Browser
$data = @(
    @("QA", "TryDermaFlash",     "Core",                   "Kit", "Chrome"),  # row 2
    @("QA", "Dr.Denese",         "Core",                   "Kit", "Chrome"),  # row 3
    @("QA", "Dr.Denese",         "wb50poff",               "Kit", "Chrome"),  # row 4
    @("QA", "PrincipalSecret",   "Core",                   "Kit", "Chrome"),  # row 5
    @("QA", "PrincipalSecret",   "Order30",                "Kit", "Chrome"),  # row 6
    @("QA", "ReclaimBotanical",  "Core",                   "Kit", "Chrome"),  # row 7
    @("QA", "SheerCover",        "Core",                   "Kit", "Chrome"),  # row 8
    @("QA", "Smileactives",      "Core",                   "Kit", "Chrome"),  # row 9
    @("QA", "Smileactives",      "core2",                  "Kit", "Chrome"),  # row 10
    @("QA", "Smileactives",      "10offdeluxe",            "Kit", "Chrome"),  # row 11
    @("QA", "Smileactives",      "specialoffer",           "Kit", "Chrome"),  # row 12
    @("QA", "Smileactives",      "sawb19",                 "Kit", "Chrome"),  # row 13
    @("QA", "SpecificBeauty",    "Core",                   "Kit", "Chrome"),  # row 14
    @("QA", "SpecificBeauty",    "deluxe-offer",           "Kit", "Chrome"),  # row 15
    @("QA", "Sub-D",             "Core",                   "Kit", "Chrome"),  # row 16
    @("QA", "Sub-D",             "cpcb2017",               "Kit", "Chrome"),  # row 17
    @("QA", "Sub-D",             "deluxe25off",            "Kit", "Chrome"),  # row 18
    @("QA", "Sub-D",             "deluxe25offp",           "Kit", "Chrome"),  # row 19
    @("QA", "Sub-D",             "cpwbunusedbdbj",         "Kit", "Chrome"),  # row 20
    @("QA", "Sub-D",             "deluxe25offp-holiday",   "Kit", "Chrome"),  # row 21
    @("QA", "TryDermaFlash",     "Core",                   "Kit", "Chrome"),  # row 22
    @("QA", "TryDermaFlash",     "pnln",                   "Kit", "Chrome"),  # row 23
    @("QA", "TryDermaFlash",     "trydermaflash",          "Kit", "Chrome")   # row 24
)

# Row 2 already has A2/C2/D2/E2 populated; only the Brand (B2) value actually changes.
$ws.Cells.Item(2, 2).Value = $data[0][1]

# Fill in new rows 3-24
for ($i = 1; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Move the trailing "End" marker from row 3 down to row 25
$ws.Cells.Item(25, 1).Value = "End"

# Update the selection to match the saved view state
$ws.Range("A21").Select()
